$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2019-12-31 00:00:00"
$ws.Range("O2").Value = 1315596935.77
$ws.Range("P2").Value = 69788551.73999999
$ws.Range("Q2").Value = 348194115.47
$ws.Range("R2").Value = 136.6289980912
$ws.Range("S2").Value = 376160289.85
$ws.Range("T2").Value = -0.6077474298
$ws.Range("U2").Value = 115010497.08
$ws.Range("V2").Value = 8.6576223693
$ws.Range("W2").Value = 660797813.95
$ws.Range("X2").Value = 282990810.53
$ws.Range("Y2").Value = 23.089919003
$ws.Range("Z2").Value = 51456649.67
$ws.Range("AA2").Value = 216.6178401985
$ws.Range("AB2").Value = 654799121.8200001
$ws.Range("AC2").Value = 35.3420571616
$ws.Range("AD2").Value = 32.4098540751
$ws.Range("AE2").Value = 29.6269636142
$ws.Range("AF2").Value = 169.1677166675
$ws.Range("AG2").Value = 50.2279836615
